$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E20").Value = 27.6
$ws.Range("E21").Value = 5
$ws.Range("E22").Value = 5
$ws.Range("E23").Value = 7.6
$ws.Range("E24").Value = 3.48
$ws.Range("B2:G26").NumberFormat = "#,##0.0"
